# Horarios actualizados Linea 141 - 426
# Scrape refresh: "Ultima actualizacion" 06:52:38 -> 07:14:27, with new rows
# of arrival data appended/inserted chronologically across the 3 sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 07:14:27"
$ws1.Range("A3").Value = "Total filas: 78"

# New row before (old) row 52
$ws1.Rows.Item(52).Insert()
$ws1.Range("A52").Value = "07:14:27"
$ws1.Range("B52").Value = "07:14"
$ws1.Range("C52").Value = "11_ETCHEVERRY"
$ws1.Range("D52").Value = 0
$ws1.Range("E52").Value = "LP1912"

# New row before (old) row 64 -> ends up at row 65
$ws1.Rows.Item(65).Insert()
$ws1.Range("A65").Value = "07:14:27"
$ws1.Range("B65").Value = "07:37"
$ws1.Range("C65").Value = "23_HERNANDEZ"
$ws1.Range("D65").Value = 23
$ws1.Range("E65").Value = "LP1912"

# Two new rows before (old) row 69 -> end up at rows 71 and 72
$ws1.Rows.Item(71).Insert()
$ws1.Range("A71").Value = "07:14:27"
$ws1.Range("B71").Value = "07:58"
$ws1.Range("C71").Value = "16_SANTA ANA"
$ws1.Range("D71").Value = 44
$ws1.Range("E71").Value = "LP1912"

$ws1.Rows.Item(72).Insert()
$ws1.Range("A72").Value = "07:14:27"
$ws1.Range("B72").Value = "08:03"
$ws1.Range("C72").Value = "11_ETCHEVERRY"
$ws1.Range("D72").Value = 49
$ws1.Range("E72").Value = "LP1912"

# Five new rows appended at the end (rows 79-83)
$ws1.Range("A79").Value = "07:14:27"
$ws1.Range("B79").Value = "08:43"
$ws1.Range("C79").Value = "14_ABASTO"
$ws1.Range("D79").Value = 89
$ws1.Range("E79").Value = "LP1912"

$ws1.Range("A80").Value = "07:14:27"
$ws1.Range("B80").Value = "08:54"
$ws1.Range("C80").Value = "17_ROMERO"
$ws1.Range("D80").Value = 100
$ws1.Range("E80").Value = "LP1912"

$ws1.Range("A81").Value = "07:14:27"
$ws1.Range("B81").Value = "09:01"
$ws1.Range("C81").Value = "215A_EL PATO"
$ws1.Range("D81").Value = 107
$ws1.Range("E81").Value = "LP1912"

$ws1.Range("A82").Value = "07:14:27"
$ws1.Range("B82").Value = "09:07"
$ws1.Range("C82").Value = "23_HERNANDEZ"
$ws1.Range("D82").Value = 113
$ws1.Range("E82").Value = "LP1912"

$ws1.Range("A83").Value = "07:14:27"
$ws1.Range("B83").Value = "09:10"
$ws1.Range("C83").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D83").Value = 116
$ws1.Range("E83").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 07:14:27"
$ws2.Range("A3").Value = "Total filas: 14"

# One new row appended at the end (row 19)
$ws2.Range("A19").Value = "07:14:27"
$ws2.Range("B19").Value = "09:01"
$ws2.Range("C19").Value = "215A_EL PATO"
$ws2.Range("D19").Value = 107
$ws2.Range("E19").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 07:14:27"
$ws3.Range("A3").Value = "Total filas: 16"

# New row before (old) row 17
$ws3.Rows.Item(17).Insert()
$ws3.Range("A17").Value = "07:14:27"
$ws3.Range("B17").Value = "08:14"
$ws3.Range("C17").Value = "215C_LA PLATA"
$ws3.Range("D17").Value = 60
$ws3.Range("E17").Value = "L6203"

# Two new rows appended at the end (rows 20 and 21)
$ws3.Range("A20").Value = "07:14:27"
$ws3.Range("B20").Value = "08:35"
$ws3.Range("C20").Value = "215A_LA PLATA"
$ws3.Range("D20").Value = 81
$ws3.Range("E20").Value = "L6173"

$ws3.Range("A21").Value = "07:14:27"
$ws3.Range("B21").Value = "09:09"
$ws3.Range("C21").Value = "215D_LA PLATA"
$ws3.Range("D21").Value = 115
$ws3.Range("E21").Value = "L6203"

Write-Host "Edit applied."
